# "Fruta / hortaliza, semanal" — add this week's new price observation.
#
# The sheet is a flat log of weekly Mango price reports (Feria Lagunitas de
# Puerto Montt). A new weekly record is inserted at row 240, which pushes
# every existing record from row 240 onward down by one row (240->241,
# ..., 357->358), growing the used range from A1:T357 to A1:T358.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 240, shifting rows 240:357 down to 241:358.
$ws.Rows(240).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A240").Value = 4
$ws.Range("B240").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C240").Value = 'Los Lagos'
$ws.Range("D240").Value = 45029
$ws.Range("E240").Value = 10
$ws.Range("F240").Value = 'Fruta'
$ws.Range("G240").Value = 100108
$ws.Range("H240").Value = 'Tropicales y subtropicales'
$ws.Range("I240").Value = 100108002
$ws.Range("J240").Value = 'Mango'
$ws.Range("K240").Value = 'Sin especificar'
$ws.Range("L240").Value = 'Primera'
$ws.Range("M240").Value = 200
$ws.Range("N240").Value = 8000
$ws.Range("O240").Value = 8500
$ws.Range("P240").Value = 8250
$ws.Range("Q240").Value = '$/bandeja 4 kilos'
$ws.Range("R240").Value = 'Perú'
$ws.Range("S240").Value = 2062
$ws.Range("T240").Value = 4
